$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column ("TabName" / "CasesTab") - shifts existing
# query/dbExcel/WebExcel columns one to the right.
$ws.Columns("A:A").Insert() | Out-Null

$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# New "StatQuery" Cypher query (now in column C, was column B).
$ws.Range("C2").Value = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)
    WHERE c.race = ""AMERICAN_INDIAN_OR_ALASKA_NATIVE""
OPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)
RETURN 
    COUNT(DISTINCT f) AS number_of_files,
    COUNT(DISTINCT c.case_id) AS number_of_cases,
    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"

# New "query" Cypher query (now in column B, was column A).
$ws.Range("B2").Value = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)
    WHERE c.race = ""AMERICAN_INDIAN_OR_ALASKA_NATIVE""
WITH DISTINCT c, a, ct
RETURN 
    COALESCE(c.case_id, '') AS ``Case ID``,
    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,
    COALESCE(a.arm_id, '') AS ``Arm``,
    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,
    COALESCE(c.disease, '') AS ``Diagnosis``,
    COALESCE(c.gender, '') AS ``Gender``,
    COALESCE(c.race, '') AS ``Race``,
    COALESCE(c.ethnicity, '') AS ``Ethnicity``"

# New narrow "best fit" width for the inserted TabName column, and a taller
# row to accommodate the longer wrapped query text in B2:C2.
$ws.Columns("A:A").ColumnWidth = 8
$ws.Rows(2).RowHeight = 174

$ws.Range("B4").Select() | Out-Null
